# Add ability to get data and send in sms.
# Inserts a new "subject_name" text question above the existing "send_sms"
# row on the survey sheet, and makes the survey sheet the active tab.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Insert a new row 2 on the survey sheet (pushes the send_sms row down to row 3)
# and populate it with the new subject_name text question.
$survey.Rows.Item(2).Insert()

$survey.Range("B2").Value = "subject_name"
$survey.Range("A2").Value = "text"
$survey.Range("C2").Value = "Enter the subject's name."
$survey.Rows.Item(2).RowHeight = 12

# Make the survey sheet the active tab/sheet (it was previously "settings"),
# and move the selection to D3, where the send_sms row's table_id value lives
# now that it has shifted down a row.
$survey.Activate()
[void]$survey.Range("D3").Select()
